$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "expectedValue"
$ws.Range("B2").Value = "bala"

$ws.Range("A4:XFD1048576").Select()
